# Compliance bugfixes: flip specific rows in column C ("compliance") from 0 to 1.
# The row numbers below correspond to the exact set of cells whose value
# changed from 0 to 1 in the source diff (sheet "Sheet1", column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(4,5,6,16,20,23,24,25,26,32,34,42,44,47,49,50,65,66,68,73,78,81,85,89,93,94,99,111,119,133,136,142,147,148,149,154,159,163,172,176,182,185,188,192,193,200,210,211,219,229,240,245,249,255,266,267,273,291,292,308,309,339,355,367,368,379,387,405,408,431,445,457,467,468,472,481,488,490,505,519,547,553,569,593,608,684,720,734,736,737,744,749,762,766,769,772,779,780,781,783,791,792,794,800,801,811,820,821,842,855,856,859,870,881,899,909,924,944,947,954,957,962,966,972,982,992,1014,1023,1027,1028,1038,1049,1068,1071,1098,1100,1110,1119,1125,1129,1174,1224,1229,1231,1252,1303,1321,1322,1332,1343,1350,1353,1362,1372,1378,1379,1384,1386,1401,1415,1438,1473,1486,1487,1488,1501,1510,1514,1535,1560,1566,1583,1586,1592,1598,1599)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 1
}
